# Add daily power records: fill in the missing Start/End Time for row 128
# (2018-12-16) and append three new daily records (2018-12-17, 2018-12-18,
# 2018-12-19) to the "comforter_cda_table" table, expanding it from
# A1:F128 to A1:F131.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 128 (2018-12-16) already exists but was missing Start/End Time values.
$ws.Range("B128").Value = 0
$ws.Range("C128").Value = 0

# Grow the table by three rows (2018-12-17, 2018-12-18, 2018-12-19) so the
# table/autofilter range and worksheet dimension expand to F131.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Row 129: 2018-12-17
$ws.Range("A129").Value = 43451
$ws.Range("B129").Value = 0
$ws.Range("C129").Value = 0
$ws.Range("D129").Formula = "=(C129-B129)* 1440"
$ws.Range("E129").Formula = "=IF(C129>B129, (C129-B129)*1440, (B129-C129)*1440)"
$ws.Range("F129").Formula = "=ABS((C129-B129)*1440)"

# Row 130: 2018-12-18
$ws.Range("A130").Value = 43452
$ws.Range("B130").Value = 0.76527777777777783
$ws.Range("C130").Value = 0.99930555555555556
$ws.Range("D130").Formula = "=(C130-B130)* 1440"
$ws.Range("E130").Formula = "=IF(C130>B130, (C130-B130)*1440, (B130-C130)*1440)"
$ws.Range("F130").Formula = "=ABS((C130-B130)*1440)"

# Row 131: 2018-12-19
$ws.Range("A131").Value = 43453
$ws.Range("B131").Value = 0
$ws.Range("D131").Formula = "=(C131-B131)* 1440"
$ws.Range("E131").Formula = "=IF(C131>B131, (C131-B131)*1440, (B131-C131)*1440)"
$ws.Range("F131").Formula = "=ABS((C131-B131)*1440)"

# Match the saved selection / scroll position from the source edit.
$ws.Range("C131").Select() | Out-Null
